# Students Detail.xlsx edit:
#  - Fix the typo in the second worksheet's name ("shet61" -> "sheet61")
#  - Make that sheet the active tab (it was previously "data" / sheet 1
#    that had the tab selected; tab selection now moves to "sheet61")

$wb = $excel.ActiveWorkbook

# The workbook has two sheets: "data" (first) and the misspelled
# "shet61" (second). Grab the second sheet and correct its name.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "sheet61"

# Select/activate the corrected sheet so it becomes the workbook's
# active tab (this also clears tabSelected on the previously-active
# "data" sheet).
$ws2.Activate()
